# Slide 1, "Subtítulo 2" placeholder (subtitle) holds two paragraphs:
#   "Nádio dib – eng. Software"
#   "uniprojeção"
# Both runs get an explicit yellow (FFFF00) font color.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)

# Setting the color on the whole TextRange applies it to every run in the
# text frame (both paragraphs), matching <a:solidFill><a:srgbClr val="FFFF00"/>
# being added to each run's rPr.
$sh.TextFrame.TextRange.Font.Color.RGB = 65535
